$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 10:52:54.185365",
    "2021-10-05 10:52:54.185376",
    "2021-10-05 10:52:54.185379",
    "2021-10-05 10:52:54.185382",
    "2021-10-05 10:52:54.185385",
    "2021-10-05 10:52:54.185387",
    "2021-10-05 10:52:54.185390",
    "2021-10-05 10:52:54.185392",
    "2021-10-05 10:52:54.185395",
    "2021-10-05 10:52:54.185398",
    "2021-10-05 10:52:54.185400",
    "2021-10-05 10:52:54.185403",
    "2021-10-05 10:52:54.185405",
    "2021-10-05 10:52:54.185407",
    "2021-10-05 10:52:54.185410",
    "2021-10-05 10:52:54.185412",
    "2021-10-05 10:52:54.185415",
    "2021-10-05 10:52:54.185418",
    "2021-10-05 10:52:54.185421",
    "2021-10-05 10:52:54.185423",
    "2021-10-05 10:52:54.185426",
    "2021-10-05 10:52:54.185428",
    "2021-10-05 10:52:54.185430",
    "2021-10-05 10:52:54.185433",
    "2021-10-05 10:52:54.185436",
    "2021-10-05 10:52:54.185438",
    "2021-10-05 10:52:54.185441",
    "2021-10-05 10:52:54.185443",
    "2021-10-05 10:52:54.185446",
    "2021-10-05 10:52:54.185448",
    "2021-10-05 10:52:54.185451",
    "2021-10-05 10:52:54.185453",
    "2021-10-05 10:52:54.185456",
    "2021-10-05 10:52:54.185459",
    "2021-10-05 10:52:54.185462",
    "2021-10-05 10:52:54.185465",
    "2021-10-05 10:52:54.185467",
    "2021-10-05 10:52:54.185470",
    "2021-10-05 10:52:54.185472",
    "2021-10-05 10:52:54.185475"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
